{"js": "const body = context.document.body;\nconst pairs = [\n  [\"167\u00d77=1169\", \"473\u00d79=4257\"],\n  [\"850\u00d72=1700\", \"944\u00d77=6608\"],\n  [\"978\u00d74=3912\", \"251\u00d77=1757\"],\n  [\"370\u00d78=2960\", \"556\u00d76=3336\"],\n  [\"357\u00d73=1071\", \"432\u00d75=2160\"],\n  [\"836\u00d76=5016\", \"283\u00d78=2264\"],\n  [\"285\u00d73=855\", \"601\u00d77=4207\"],\n  [\"476\u00d76=2856\", \"418\u00d78=3344\"],\n  [\"399\u00d76=2394\", \"788\u00d75=3940\"],\n  [\"701\u00d74=2804\", \"638\u00d72=1276\"],\n  [\"242\u00d77=1694\", \"511\u00d77=3577\"],\n  [\"634\u00d76=3804\", \"353\u00d77=2471\"],\n  [\"201\u00d76=1206\", \"323\u00d73=969\"],\n  [\"119\u00d75=595\", \"598\u00d73=1794\"],\n  [\"209\u00d72=418\", \"647\u00d78=5176\"],\n  [\"815\u00d72=1630\", \"531\u00d75=2655\"],\n  [\"258\u00d75=1290\", \"438\u00d75=2190\"],\n  [\"237\u00d72=474\", \"146\u00d76=876\"],\n  [\"814\u00d77=5698\", \"366\u00d74=1464\"],\n  [\"868\u00d74=3472\", \"126\u00d72=252\"],\n  [\"237\u00d74=948\", \"589\u00d73=1767\"],\n  [\"587\u00d78=4696\", \"597\u00d79=5373\"],\n  [\"427\u00d75=2135\", \"267\u00d78=2136\"],\n  [\"991\u00d73=2973\", \"432\u00d77=3024\"],\n  [\"598\u00d79=5382\", \"582\u00d73=1746\"]\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n$d.Content.Find.Execute(\"167\u00d77=1169\", $false, $false, $false, $false, $false, $true, 1, $false, \"473\u00d79=4257\", 2) | Out-Null\n$d.Content.Find.Execute(\"850\u00d72=1700\", $false, $false, $false, $false, $false, $true, 1, $false, \"944\u00d77=6608\", 2) | Out-Null\n$d.Content.Find.Execute(\"978\u00d74=3912\", $false, $false, $false, $false, $false, $true, 1, $false, \"251\u00d77=1757\", 2) | Out-Null\n$d.Content.Find.Execute(\"370\u00d78=2960\", $false, $false, $false, $false, $false, $true, 1, $false, \"556\u00d76=3336\", 2) | Out-Null\n$d.Content.Find.Execute(\"357\u00d73=1071\", $false, $false, $false, $false, $false, $true, 1, $false, \"432\u00d75=2160\", 2) | Out-Null\n$d.Content.Find.Execute(\"836\u00d76=5016\", $false, $false, $false, $false, $false, $true, 1, $false, \"283\u00d78=2264\", 2) | Out-Null\n$d.Content.Find.Execute(\"285\u00d73=855\", $false, $false, $false, $false, $false, $true, 1, $false, \"601\u00d77=4207\", 2) | Out-Null\n$d.Content.Find.Execute(\"476\u00d76=2856\", $false, $false, $false, $false, $false, $true, 1, $false, \"418\u00d78=3344\", 2) | Out-Null\n$d.Content.Find.Execute(\"399\u00d76=2394\", $false, $false, $false, $false, $false, $true, 1, $false, \"788\u00d75=3940\", 2) | Out-Null\n$d.Content.Find.Execute(\"701\u00d74=2804\", $false, $false, $false, $false, $false, $true, 1, $false, \"638\u00d72=1276\", 2) | Out-Null\n$d.Content.Find.Execute(\"242\u00d77=1694\", $false, $false, $false, $false, $false, $true, 1, $false, \"511\u00d77=3577\", 2) | Out-Null\n$d.Content.Find.Execute(\"634\u00d76=3804\", $false, $false, $false, $false, $false, $true, 1, $false, \"353\u00d77=2471\", 2) | Out-Null\n$d.Content.Find.Execute(\"201\u00d76=1206\", $false, $false, $false, $false, $false, $true, 1, $false, \"323\u00d73=969\", 2) | Out-Null\n$d.Content.Find.Execute(\"119\u00d75=595\", $false, $false, $false, $false, $false, $true, 1, $false, \"598\u00d73=1794\", 2) | Out-Null\n$d.Content.Find.Execute(\"209\u00d72=418\", $false, $false, $false, $false, $false, $true, 1, $false, \"647\u00d78=5176\", 2) | Out-Null\n$d.Content.Find.Execute(\"815\u00d72=1630\", $false, $false, $false, $false, $false, $true, 1, $false, \"531\u00d75=2655\", 2) | Out-Null\n$d.Content.Find.Execute(\"258\u00d75=1290\", $false, $false, $false, $false, $false, $true, 1, $false, \"438\u00d75=2190\", 2) | Out-Null\n$d.Content.Find.Execute(\"237\u00d72=474\", $false, $false, $false, $false, $false, $true, 1, $false, \"146\u00d76=876\", 2) | Out-Null\n$d.Content.Find.Execute(\"814\u00d77=5698\", $false, $false, $false, $false, $false, $true, 1, $false, \"366\u00d74=1464\", 2) | Out-Null\n$d.Content.Find.Execute(\"868\u00d74=3472\", $false, $false, $false, $false, $false, $true, 1, $false, \"126\u00d72=252\", 2) | Out-Null\n$d.Content.Find.Execute(\"237\u00d74=948\", $false, $false, $false, $false, $false, $true, 1, $false, \"589\u00d73=1767\", 2) | Out-Null\n$d.Content.Find.Execute(\"587\u00d78=4696\", $false, $false, $false, $false, $false, $true, 1, $false, \"597\u00d79=5373\", 2) | Out-Null\n$d.Content.Find.Execute(\"427\u00d75=2135\", $false, $false, $false, $false, $false, $true, 1, $false, \"267\u00d78=2136\", 2) | Out-Null\n$d.Content.Find.Execute(\"991\u00d73=2973\", $false, $false, $false, $false, $false, $true, 1, $false, \"432\u00d77=3024\", 2) | Out-Null\n$d.Content.Find.Execute(\"598\u00d79=5382\", $false, $false, $false, $false, $false, $true, 1, $false, \"582\u00d73=1746\", 2) | Out-Null\n"}
